$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a union of all target cells so the Text number format can be
# applied in one pass (keeps the style table minimal/deduped), then
# assign each literal text value individually. Union.NumberFormat only
# reliably reaches the first Area, so iterate over .Areas explicitly.
$targets = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11",
    "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16",
    "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22",
    "E22", "D23", "E23", "D24", "E24", "E25", "E26", "D27", "E27", "D39",
    "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "E44",
    "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49",
    "D50", "E50", "D51", "E51"
)

$union = $null
foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    if ($union -eq $null) { $union = $cell } else { $union = $excel.Union($union, $cell) }
}

# Format every touched cell as Text first so that numeric-looking
# literals (e.g. "304.53", "1.10%") are stored verbatim instead of being
# reinterpreted by Excel as numbers/percentages.
foreach ($area in $union.Areas) {
    $area.NumberFormat = "@"
}

# Now write each new literal value.
$ws.Range("D2").Value = "304.53"
$ws.Range("E2").Value = "1.10%"
$ws.Range("D3").Value = "35.60"
$ws.Range("E3").Value = "1.33%"
$ws.Range("D4").Value = "5.061"
$ws.Range("E4").Value = "0.48%"
$ws.Range("D5").Value = "0.08047"
$ws.Range("E5").Value = "0.82%"
$ws.Range("D6").Value = "1.911"
$ws.Range("E6").Value = "0.28%"
$ws.Range("D7").Value = "4.179"
$ws.Range("E7").Value = "3.15%"
$ws.Range("D8").Value = "7.737"
$ws.Range("E8").Value = "-0.84%"
$ws.Range("D9").Value = "0.9282"
$ws.Range("E9").Value = "0.63%"
$ws.Range("D10").Value = "0.1385"
$ws.Range("E10").Value = "6.68%"
$ws.Range("D11").Value = "0.1898"
$ws.Range("E11").Value = "2.37%"
$ws.Range("D12").Value = "0.09136"
$ws.Range("E12").Value = "-7.62%"
$ws.Range("D13").Value = "0.03631"
$ws.Range("E13").Value = "2.65%"
$ws.Range("D14").Value = "0.09812"
$ws.Range("E14").Value = "-0.41%"
$ws.Range("D15").Value = "0.001433"
$ws.Range("E15").Value = "3.17%"
$ws.Range("D16").Value = "0.005907"
$ws.Range("E16").Value = "0.56%"
$ws.Range("D17").Value = "3.554"
$ws.Range("E17").Value = "1.37%"
$ws.Range("E18").Value = "-1.35%"
$ws.Range("D19").Value = "0.3454"
$ws.Range("E19").Value = "1.58%"
$ws.Range("D20").Value = "0.1331"
$ws.Range("E20").Value = "2.35%"
$ws.Range("D21").Value = "4.890"
$ws.Range("E21").Value = "-3.02%"
$ws.Range("D22").Value = "0.2510"
$ws.Range("E22").Value = "4.57%"
$ws.Range("D23").Value = "0.04437"
$ws.Range("E23").Value = "-1.33%"
$ws.Range("D24").Value = "0.001222"
$ws.Range("E24").Value = "0.47%"
$ws.Range("E25").Value = "-0.15%"
$ws.Range("E26").Value = "24.82%"
$ws.Range("D27").Value = "0.0003133"
$ws.Range("E27").Value = "4.32%"
$ws.Range("D39").Value = "0.01959"
$ws.Range("E39").Value = "4.13%"
$ws.Range("D40").Value = "0.04884"
$ws.Range("E40").Value = "3.50%"
$ws.Range("D41").Value = "0.007639"
$ws.Range("E41").Value = "1.26%"
$ws.Range("D42").Value = "0.009197"
$ws.Range("E42").Value = "-9.97%"
$ws.Range("D43").Value = "0.1372"
$ws.Range("E43").Value = "3.57%"
$ws.Range("E44").Value = "-0.46%"
$ws.Range("D45").Value = "0.01135"
$ws.Range("E45").Value = "7.36%"
$ws.Range("D46").Value = "0.00006406"
$ws.Range("E46").Value = "2.63%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.01%"
$ws.Range("D48").Value = "64.67"
$ws.Range("E48").Value = "0.29%"
$ws.Range("D49").Value = "0.001192"
$ws.Range("E49").Value = "-19.96%"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "0.01%"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "0.01%"
